# Updates cryptos list values (Price / Volume(1h) / Coin / Link columns)
# to match the latest scrape, per the authoritative diff.
#
# Notes:
#  - The Price/Volume cells are stored as literal text (t="inlineStr" in the
#    original OOXML), not numbers, even though many of them look numeric
#    (e.g. "1.000", "0.9995"). Excel's COM layer auto-converts a plain
#    assignment of such a numeric-looking string into a true number, which
#    would change the cell type. To keep these as text we temporarily force
#    the cell's number format to Text ("@") before assigning the value, then
#    restore the cell style to "Normal" afterwards so no stray formatting
#    is left behind on the cell.
#  - A couple of row pairs (28/29, 48/49, and 50/51 shifting into a new 50)
#    had their Coin name + Link swapped/replaced along with Price/Volume, so
#    columns B and C are updated for those rows too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

function Set-PlainValue($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

Set-PlainValue "D2" "31.219.62"
Set-PlainValue "E2" "  +2.05%  "
Set-PlainValue "D3" "1.995.00"
Set-PlainValue "E3" "  +5.90%  "
Set-TextValue "D4" "1.000"
Set-PlainValue "E4" "  +0.00%  "
Set-TextValue "D5" "0.7769"
Set-PlainValue "E5" "  +63.88%  "
Set-TextValue "D6" "254.12"
Set-PlainValue "E6" "  +3.22%  "
Set-TextValue "D7" "0.9995"
Set-PlainValue "E7" "  -0.08%  "
Set-TextValue "D8" "0.3476"
Set-PlainValue "E8" "  +20.06%  "
Set-TextValue "D9" "27.83"
Set-PlainValue "E9" "  +24.51%  "
Set-TextValue "D10" "0.07062"
Set-TextValue "D11" "0.8397"
Set-PlainValue "E11" "  +10.08%  "
Set-TextValue "D12" "0.08199"
Set-PlainValue "E12" "  +4.80%  "
Set-TextValue "D13" "100.70"
Set-PlainValue "E13" "  +1.05%  "
Set-PlainValue "D14" "1.995.31"
Set-PlainValue "E14" "  +5.92%  "
Set-TextValue "D15" "5.624"
Set-PlainValue "E15" "  +7.45%  "
Set-TextValue "D16" "15.16"
Set-PlainValue "E16" "  +15.07%  "
Set-TextValue "D17" "272.11"
Set-PlainValue "E17" "  -4.11%  "
Set-PlainValue "D18" "31.221.18"
Set-PlainValue "E18" "  +2.13%  "
Set-TextValue "D19" "5.960"
Set-PlainValue "E19" "  +11.37%  "
Set-TextValue "D20" "0.000008002"
Set-PlainValue "E20" "  +6.37%  "
Set-PlainValue "D21" "2.258.58"
Set-PlainValue "E21" "  +6.35%  "
Set-TextValue "D22" "0.9991"
Set-PlainValue "E22" "  -0.06%  "
Set-TextValue "D23" "0.9991"
Set-PlainValue "E23" "  -0.12%  "
Set-TextValue "D24" "7.086"
Set-PlainValue "E24" "  +10.14%  "
Set-TextValue "D25" "9.979"
Set-PlainValue "E25" "  +8.86%  "
Set-TextValue "D26" "164.76"
Set-PlainValue "E26" "  +0.56%  "
Set-TextValue "D27" "0.1430"
Set-PlainValue "E27" "  +46.62%  "
Set-PlainValue "B28" "LidoDAOToken"
Set-PlainValue "C28" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D28" "2.419"
Set-PlainValue "E28" "  +26.91%  "
Set-PlainValue "B29" "EthereumClassic"
Set-PlainValue "C29" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D29" "19.87"
Set-PlainValue "E29" "  +4.53%  "
Set-TextValue "D30" "1.593"
Set-PlainValue "E30" "  +6.02%  "
Set-TextValue "D31" "1.363"
Set-PlainValue "E31" "  +2.66%  "
Set-TextValue "D32" "4.601"
Set-PlainValue "E32" "  +8.27%  "
Set-TextValue "D33" "4.442"
Set-PlainValue "E33" "  +6.13%  "
Set-TextValue "D34" "0.05327"
Set-PlainValue "E34" "  +10.01%  "
Set-TextValue "D35" "1.245"
Set-PlainValue "E35" "  +9.96%  "
Set-TextValue "D36" "0.7913"
Set-PlainValue "E36" "  +13.32%  "
Set-TextValue "D37" "2.766"
Set-PlainValue "E37" "  -0.24%  "
Set-TextValue "D38" "0.9986"
Set-PlainValue "E38" "  -0.12%  "
Set-TextValue "D39" "0.02002"
Set-PlainValue "E39" "  +5.13%  "
Set-PlainValue "E40" "  +1.50%  "
Set-TextValue "D41" "83.25"
Set-PlainValue "E41" "  +10.28%  "
Set-TextValue "D42" "6.773"
Set-PlainValue "E42" "  +7.42%  "
Set-TextValue "D43" "0.4670"
Set-PlainValue "E43" "  +9.88%  "
Set-TextValue "D44" "2.132"
Set-PlainValue "E44" "  +7.93%  "
Set-TextValue "D45" "0.8559"
Set-PlainValue "E45" "  +2.09%  "
Set-TextValue "D46" "105.09"
Set-PlainValue "E46" "  +3.63%  "
Set-TextValue "D47" "0.9993"
Set-PlainValue "E47" "  -0.09%  "
Set-PlainValue "B48" "Aptos"
Set-PlainValue "C48" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D48" "7.716"
Set-PlainValue "E48" "  +9.95%  "
Set-PlainValue "B49" "EnergySwap"
Set-PlainValue "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "10.01"
Set-PlainValue "E49" "  +0.46%  "
Set-PlainValue "B50" "SynthetixNetwork"
Set-PlainValue "C50" "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextValue "D50" "3.038"
Set-PlainValue "E50" "  +47.06%  "
Set-PlainValue "B51" "Elrond"
Set-PlainValue "C51" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue "D51" "37.56"
Set-PlainValue "E51" "  +6.53%  "
